# Apply the "update scripts with new tpm" edit to the NATMI LR-pair sheet.
# The underlying TPM recalculation changed Ligand-expressing-cell counts and all
# downstream statistics for every row, and the "Resolving-Mac" sending-cluster block
# (rows 14:17) is dropped entirely now that Myoc is no longer detected as sent from it.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the four rows for Sending cluster = Resolving-Mac (rows 14-17 in the original sheet).
$ws.Range("A14:T17").EntireRow.Delete()

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Myoc"
$row2[0,2] = "Fzd7"
$row2[0,3] = "ECs"
$row2[0,4] = 1
$row2[0,5] = 0.3333333333333333
$row2[0,6] = 0.2439023333333333
$row2[0,7] = 0.731707
$row2[0,8] = 0.006368708152767561
$row2[0,9] = 0.006368708152767561
$row2[0,10] = 3
$row2[0,11] = 1
$row2[0,12] = 1.123319
$row2[0,13] = 3.369957
$row2[0,14] = 0.05053686506648315
$row2[0,15] = 0.05053686506648315
$row2[0,16] = 0.2739801251776667
$row2[0,17] = 2.465821126599
$row2[0,18] = 0.0003218545445642254
$row2[0,19] = 0.0003218545445642254
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Myoc"
$row3[0,2] = "Fzd7"
$row3[0,3] = "FAPs"
$row3[0,4] = 1
$row3[0,5] = 0.3333333333333333
$row3[0,6] = 0.2439023333333333
$row3[0,7] = 0.731707
$row3[0,8] = 0.006368708152767561
$row3[0,9] = 0.006368708152767561
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 11.20764866666667
$row3[0,13] = 33.622946
$row3[0,14] = 0.5042195746532222
$row3[0,15] = 0.5042195746532223
$row3[0,16] = 2.733571660980222
$row3[0,17] = 24.602144948822
$row3[0,18] = 0.003211227315878968
$row3[0,19] = 0.003211227315878969
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Myoc"
$row4[0,2] = "Fzd7"
$row4[0,3] = "MuSCs"
$row4[0,4] = 1
$row4[0,5] = 0.3333333333333333
$row4[0,6] = 0.2439023333333333
$row4[0,7] = 0.731707
$row4[0,8] = 0.006368708152767561
$row4[0,9] = 0.006368708152767561
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 4.958620666666667
$row4[0,13] = 14.875862
$row4[0,14] = 0.2230827962023326
$row4[0,15] = 0.2230827962023326
$row4[0,16] = 1.209419150714889
$row4[0,17] = 10.884772356434
$row4[0,18] = 0.00142074922291598
$row4[0,19] = 0.00142074922291598
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "ECs"
$row5[0,1] = "Myoc"
$row5[0,2] = "Fzd7"
$row5[0,3] = "Resolving-Mac"
$row5[0,4] = 1
$row5[0,5] = 0.3333333333333333
$row5[0,6] = 0.2439023333333333
$row5[0,7] = 0.731707
$row5[0,8] = 0.006368708152767561
$row5[0,9] = 0.006368708152767561
$row5[0,10] = 3
$row5[0,11] = 1
$row5[0,12] = 4.938126
$row5[0,13] = 14.814378
$row5[0,14] = 0.222160764077962
$row5[0,15] = 0.222160764077962
$row5[0,16] = 1.204420453694
$row5[0,17] = 10.839784083246
$row5[0,18] = 0.001414877069408387
$row5[0,19] = 0.001414877069408387
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "FAPs"
$row6[0,1] = "Myoc"
$row6[0,2] = "Fzd7"
$row6[0,3] = "ECs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 37.428193
$row6[0,7] = 112.284579
$row6[0,8] = 0.9773142989029397
$row6[0,9] = 0.9773142989029399
$row6[0,10] = 3
$row6[0,11] = 1
$row6[0,12] = 1.123319
$row6[0,13] = 3.369957
$row6[0,14] = 0.05053686506648315
$row6[0,15] = 0.05053686506648315
$row6[0,16] = 42.04380033256701
$row6[0,17] = 378.3942029931031
$row6[0,18] = 0.04939040085120244
$row6[0,19] = 0.04939040085120245
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "FAPs"
$row7[0,1] = "Myoc"
$row7[0,2] = "Fzd7"
$row7[0,3] = "FAPs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 37.428193
$row7[0,7] = 112.284579
$row7[0,8] = 0.9773142989029397
$row7[0,9] = 0.9773142989029399
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 11.20764866666667
$row7[0,13] = 33.622946
$row7[0,14] = 0.5042195746532222
$row7[0,15] = 0.5042195746532223
$row7[0,16] = 419.4820373721926
$row7[0,17] = 3775.338336349734
$row7[0,18] = 0.4927810000953524
$row7[0,19] = 0.4927810000953525
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "FAPs"
$row8[0,1] = "Myoc"
$row8[0,2] = "Fzd7"
$row8[0,3] = "MuSCs"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 37.428193
$row8[0,7] = 112.284579
$row8[0,8] = 0.9773142989029397
$row8[0,9] = 0.9773142989029399
$row8[0,10] = 3
$row8[0,11] = 1
$row8[0,12] = 4.958620666666667
$row8[0,13] = 14.875862
$row8[0,14] = 0.2230827962023326
$row8[0,15] = 0.2230827962023326
$row8[0,16] = 185.5922113257887
$row8[0,17] = 1670.329901932098
$row8[0,18] = 0.2180220065677901
$row8[0,19] = 0.2180220065677901
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "FAPs"
$row9[0,1] = "Myoc"
$row9[0,2] = "Fzd7"
$row9[0,3] = "Resolving-Mac"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 37.428193
$row9[0,7] = 112.284579
$row9[0,8] = 0.9773142989029397
$row9[0,9] = 0.9773142989029399
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 4.938126
$row9[0,13] = 14.814378
$row9[0,14] = 0.222160764077962
$row9[0,15] = 0.222160764077962
$row9[0,16] = 184.825132986318
$row9[0,17] = 1663.426196876862
$row9[0,18] = 0.2171208913885948
$row9[0,19] = 0.2171208913885948
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = "MuSCs"
$row10[0,1] = "Myoc"
$row10[0,2] = "Fzd7"
$row10[0,3] = "ECs"
$row10[0,4] = 3
$row10[0,5] = 1
$row10[0,6] = 0.6248916666666666
$row10[0,7] = 1.874675
$row10[0,8] = 0.01631699294429263
$row10[0,9] = 0.01631699294429263
$row10[0,10] = 3
$row10[0,11] = 1
$row10[0,12] = 1.123319
$row10[0,13] = 3.369957
$row10[0,14] = 0.05053686506648315
$row10[0,15] = 0.05053686506648315
$row10[0,16] = 0.7019526821083334
$row10[0,17] = 6.317574138975
$row10[0,18] = 0.0008246096707164741
$row10[0,19] = 0.0008246096707164741
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,20
$row11[0,0] = "MuSCs"
$row11[0,1] = "Myoc"
$row11[0,2] = "Fzd7"
$row11[0,3] = "FAPs"
$row11[0,4] = 3
$row11[0,5] = 1
$row11[0,6] = 0.6248916666666666
$row11[0,7] = 1.874675
$row11[0,8] = 0.01631699294429263
$row11[0,9] = 0.01631699294429263
$row11[0,10] = 3
$row11[0,11] = 1
$row11[0,12] = 11.20764866666667
$row11[0,13] = 33.622946
$row11[0,14] = 0.5042195746532222
$row11[0,15] = 0.5042195746532223
$row11[0,16] = 7.003566254727777
$row11[0,17] = 63.03209629254999
$row11[0,18] = 0.008227347241990857
$row11[0,19] = 0.008227347241990859
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,20
$row12[0,0] = "MuSCs"
$row12[0,1] = "Myoc"
$row12[0,2] = "Fzd7"
$row12[0,3] = "MuSCs"
$row12[0,4] = 3
$row12[0,5] = 1
$row12[0,6] = 0.6248916666666666
$row12[0,7] = 1.874675
$row12[0,8] = 0.01631699294429263
$row12[0,9] = 0.01631699294429263
$row12[0,10] = 3
$row12[0,11] = 1
$row12[0,12] = 4.958620666666667
$row12[0,13] = 14.875862
$row12[0,14] = 0.2230827962023326
$row12[0,15] = 0.2230827962023326
$row12[0,16] = 3.098600732761112
$row12[0,17] = 27.88740659485
$row12[0,18] = 0.003640040411626532
$row12[0,19] = 0.003640040411626532
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,20
$row13[0,0] = "MuSCs"
$row13[0,1] = "Myoc"
$row13[0,2] = "Fzd7"
$row13[0,3] = "Resolving-Mac"
$row13[0,4] = 3
$row13[0,5] = 1
$row13[0,6] = 0.6248916666666666
$row13[0,7] = 1.874675
$row13[0,8] = 0.01631699294429263
$row13[0,9] = 0.01631699294429263
$row13[0,10] = 3
$row13[0,11] = 1
$row13[0,12] = 4.938126
$row13[0,13] = 14.814378
$row13[0,14] = 0.222160764077962
$row13[0,15] = 0.222160764077962
$row13[0,16] = 3.08579378635
$row13[0,17] = 27.77214407715
$row13[0,18] = 0.003624995619958765
$row13[0,19] = 0.003624995619958765
$ws.Range("A13:T13").Value = $row13

# Dimension should now be A1:T13 (handled automatically by the engine based on used range).
